$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("M2").Value = 23.28617366666667
$ws.Range("N2").Value = 69.858521
$ws.Range("O2").Value = 0.2304887056246027
$ws.Range("P2").Value = 0.2304887056246027
$ws.Range("Q2").Value = 38.29938303213989
$ws.Range("R2").Value = 344.694447289259
$ws.Range("S2").Value = 0.007109102727553473
$ws.Range("T2").Value = 0.007109102727553472

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.644726333333334
$ws.Range("H3").Value = 4.934179
$ws.Range("I3").Value = 0.03084360558270512
$ws.Range("J3").Value = 0.03084360558270512
$ws.Range("O3").Value = 0.007098179626924059
$ws.Range("P3").Value = 0.007098179626924059
$ws.Range("Q3").Value = 1.179476016518
$ws.Range("R3").Value = 10.615284148662
$ws.Range("S3").Value = 0.0002189334527680387
$ws.Range("T3").Value = 0.0002189334527680387

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.644726333333334
$ws.Range("H4").Value = 4.934179
$ws.Range("I4").Value = 0.03084360558270512
$ws.Range("J4").Value = 0.03084360558270512
$ws.Range("O4").Value = 0.7624131147484733
$ws.Range("P4").Value = 0.7624131147484732
$ws.Range("Q4").Value = 126.6871269520538
$ws.Range("R4").Value = 1140.184142568484
$ws.Range("S4").Value = 0.02351556940238361
$ws.Range("T4").Value = 0.02351556940238361

$ws.Range("I5").Value = 0.828024694817689
$ws.Range("J5").Value = 0.828024694817689
$ws.Range("M5").Value = 23.28617366666667
$ws.Range("N5").Value = 69.858521
$ws.Range("O5").Value = 0.2304887056246027
$ws.Range("P5").Value = 0.2304887056246027
$ws.Range("Q5").Value = 1028.181833730739
$ws.Range("R5").Value = 9253.636503576652
$ws.Range("S5").Value = 0.1908503401337358
$ws.Range("T5").Value = 0.1908503401337358

$ws.Range("I6").Value = 0.828024694817689
$ws.Range("J6").Value = 0.828024694817689
$ws.Range("O6").Value = 0.007098179626924059
$ws.Range("P6").Value = 0.007098179626924059
$ws.Range("S6").Value = 0.005877468019344932
$ws.Range("T6").Value = 0.005877468019344932

$ws.Range("I7").Value = 0.828024694817689
$ws.Range("J7").Value = 0.828024694817689
$ws.Range("O7").Value = 0.7624131147484733
$ws.Range("P7").Value = 0.7624131147484732
$ws.Range("S7").Value = 0.6312968866646084
$ws.Range("T7").Value = 0.6312968866646083

$ws.Range("I8").Value = 0.1411316995996059
$ws.Range("J8").Value = 0.1411316995996059
$ws.Range("M8").Value = 23.28617366666667
$ws.Range("N8").Value = 69.858521
$ws.Range("O8").Value = 0.2304887056246027
$ws.Range("P8").Value = 0.2304887056246027
$ws.Range("Q8").Value = 175.2472487838157
$ws.Range("R8").Value = 1577.225239054341
$ws.Range("S8").Value = 0.03252926276331342
$ws.Range("T8").Value = 0.03252926276331342

$ws.Range("I9").Value = 0.1411316995996059
$ws.Range("J9").Value = 0.1411316995996059
$ws.Range("O9").Value = 0.007098179626924059
$ws.Range("P9").Value = 0.007098179626924059
$ws.Range("S9").Value = 0.001001778154811089
$ws.Range("T9").Value = 0.001001778154811089

$ws.Range("I10").Value = 0.1411316995996059
$ws.Range("J10").Value = 0.1411316995996059
$ws.Range("O10").Value = 0.7624131147484733
$ws.Range("P10").Value = 0.7624131147484732
$ws.Range("S10").Value = 0.1076006586814814
$ws.Range("T10").Value = 0.1076006586814814
